# Update countries & provincias Spain
# - Refresh case counts for "Estados Unidos" (row 4, unchanged ranking)
# - "Japon" jumped ahead of "Ecuador"/"Peru"/"Pakistan" in the case-count ranking:
#   remove its old row and re-insert it (with refreshed numbers) right after "Rumania"
# - "San Marino" jumped ahead of "Malta": remove its old row and re-insert it
#   (with refreshed numbers) right after "Honduras"
# - Refresh case counts for "Reunion" (row 96, unchanged ranking)
# - Refresh the "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $values) {
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
    $ws.Cells.Item($row, 7).Value = $values[6]
    $ws.Cells.Item($row, 8).Value = $values[7]
}

# Estados Unidos (row 4) - refreshed totals
Set-CountryRow 4 @("Estados Unidos", 428355, 28020, 22356, 391263, 9265, 1895, 14736)

# Japon moves from row 36 up to row 33 (right after Rumania, ahead of Ecuador/Peru/Pakistan)
$ws.Rows.Item(36).EntireRow.Delete()
$ws.Rows.Item(33).EntireRow.Insert()
Set-CountryRow 33 @("Japon", 4667, 410, 632, 3941, 99, 1, 94)

# Reunion (row 96) - refreshed totals
Set-CountryRow 96 @("Reunion", 362, 4, 40, 322, 4, 0, 0)

# San Marino moves from row 102 up to row 101 (right after Honduras, ahead of Malta)
$ws.Rows.Item(102).EntireRow.Delete()
$ws.Rows.Item(101).EntireRow.Insert()
Set-CountryRow 101 @("San Marino", 308, 29, 45, 229, 14, 0, 34)

# Updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 02:22"
